$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.990.60'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '4.032.01'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '537.21'
$c.Style = $origStyle
$ws.Range('E5').Value = '  +0.84%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '149.33'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('D7').Value = '4.027.46'
$ws.Range('E7').Value = '  -0.42%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.694'
$c.Style = $origStyle
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -1.98%  '
$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.171'
$c.Style = $origStyle
$ws.Range('E11').Value = '  -3.43%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '53.52'
$c.Style = $origStyle
$ws.Range('E12').Value = '  +6.69%  '
$c = $ws.Range('D13')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000326'
$c.Style = $origStyle
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '4.671.23'
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '4.029.06'
$ws.Range('E16').Value = '  -1.20%  '
$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.21'
$c.Style = $origStyle
$ws.Range('E17').Value = '  -1.87%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '20.76'
$c.Style = $origStyle
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('E19').Value = '  -2.87%  '
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').Value = '72.060.25'
$ws.Range('E21').Value = '  -0.05%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '432.39'
$c.Style = $origStyle
$ws.Range('E22').Value = '  -1.58%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '98.04'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -2.52%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.54'
$c.Style = $origStyle
$ws.Range('E24').Value = '  -4.94%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.26'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.69'
$c.Style = $origStyle
$ws.Range('E26').Value = '  -1.13%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.45'
$c.Style = $origStyle
$ws.Range('E27').Value = '  +27.11%  '
$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = $origStyle
$ws.Range('E28').Value = '  -1.19%  '
$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.75'
$c.Style = $origStyle
$ws.Range('E29').Value = '  -2.48%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.95'
$c.Style = $origStyle
$ws.Range('E30').Value = '  +1.70%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '37.02'
$c.Style = $origStyle
$ws.Range('E31').Value = '  -1.29%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.23'
$c.Style = $origStyle
$ws.Range('E32').Value = '  +21.62%  '
$ws.Range('E33').Value = '  +1.88%  '
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '50.34'
$c.Style = $origStyle
$ws.Range('E34').Value = '  +17.36%  '
$ws.Range('E35').Value = '  -1.06%  '
$c = $ws.Range('D36')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '680.21'
$c.Style = $origStyle
$ws.Range('E36').Value = '  +0.71%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '67.66'
$c.Style = $origStyle
$ws.Range('E37').Value = '  +1.23%  '
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.462'
$c.Style = $origStyle
$ws.Range('E38').Value = '  +5.69%  '
$ws.Range('D39').Value = '0.0₃0817'
$ws.Range('E39').Value = '  -6.12%  '
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.42'
$c.Style = $origStyle
$ws.Range('E40').Value = '  +7.35%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D41')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.148'
$c.Style = $origStyle
$ws.Range('E41').Value = '  -6.28%  '
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.38'
$c.Style = $origStyle
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.16'
$c.Style = $origStyle
$ws.Range('E43').Value = '  +16.52%  '
$ws.Range('B44').Value = 'Dai'
$ws.Range('C44').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range('E44').Value = '  -0.06%  '
$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0493'
$c.Style = $origStyle
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('E47').Value = '  -2.31%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.66'
$c.Style = $origStyle
$ws.Range('E48').Value = '  -3.69%  '
$c = $ws.Range('D49')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.35'
$c.Style = $origStyle
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '2.879.20'
$ws.Range('E51').Value = '  +10.18%  '
